$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet (becomes the new last sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "TopRelationships"

# Populate the new "TopRelationships" sheet with its label column.
$newSheet.Range("A1").Value = "Contact Name:"
$newSheet.Range("A2").Value = "Strength Rating:"
$newSheet.Range("A3").Value = "Type:"
$newSheet.Range("A4").Value = "# Activities:"
$newSheet.Range("A5").Value = "Last Activity Date:"

# Match the author's last selection/active-tab state for the new sheet.
[void]$newSheet.Range("K21").Select()
